$d = $word.ActiveDocument

# Step 1: Insert the 5 new reference paragraphs after the current last paragraph
# (the bookmarkStart/bookmarkEnd end up attached to the new last paragraph).
$insertXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Matthew Brown and David G. Lowe</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Automatic Panoramic Image Stitching using Invariant Features</w:t></w:r><w:r><w:t xml:space="preserve"> 2007</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Rolf G. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Kuehni</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Color Space and Its Divisions: Color Order from Antiquity to the Present</w:t></w:r><w:r><w:t>, 2003</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">Konstantinos N. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t>Plataniotis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t>Anastasios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> N. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t>Venetsanopoulos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">Color </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Image</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Processing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Applications</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>, 2000</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Martin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Loesdau</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>;</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Sébastien </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Chabrier</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>;</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> and Alban </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Gabillon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Hue and Saturation in the RGB Color Space</w:t></w:r><w:r><w:t>, 2014</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">ENVI - </w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">Guia do ENVI em Português. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Sulsoft</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">, 2000. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Diponível</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> em www.sulsoft.com.br. Acesso em 15 out. 2006.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endPos = $d.Content.End
$insertRange = $d.Range($endPos, $endPos)
$insertRange.InsertXML($insertXml) | Out-Null

# Step 2: Replace the "Marcio Portes..." paragraph (now paragraph 5, no longer last)
# to drop the bold paragraph-mark formatting (<w:b/> removed from pPr/rPr).
$para5Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t>Márcio Portes de Albuquerque</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t>;</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> Marcelo Portes de Albuquerque</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Processamento de Imagens: Métodos e Análises</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertXML($para5Xml) | Out-Null

Write-Host "Edit applied."
